$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()
$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 378535468.67
$ws.Range("P2").Value = 76730839.48999999
$ws.Range("Q2").Value = 17864613.67
$ws.Range("S2").Value = 105688532.98
$ws.Range("U2").Value = 61599392.4
$ws.Range("W2").Value = 200033652.99
$ws.Range("X2").Value = 48387690.64
$ws.Range("Z2").Value = 23347842.04
$ws.Range("AB2").Value = 178501815.68
$ws.Range("AF2").Value = 174.0570505753
$ws.Range("AG2").Value = 52.8440977256
